$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 306723.38
$ws.Range("J17").Value = 306723.38
$ws.Range("L17").Value = 920170.14
$ws.Range("N17").Value = -920506.14

$ws.Range("H32").Value = 2033.8334
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 2036.909
$ws.Range("K32").Value = 2000
$ws.Range("L32").Value = 2036.909
$ws.Range("M32").Value = -1674
$ws.Range("N32").Value = -2688.909

$ws.Range("H51").Value = 170293300
$ws.Range("I51").Value = 510869570
$ws.Range("K51").Value = 510869570
$ws.Range("M51").Value = -510869086

$ws.Range("H76").Value = 4999
$ws.Range("I76").Value = 4999
$ws.Range("K76").Value = 4999
$ws.Range("M76").Value = -4684

$ws.Range("H79").Value = 4999
$ws.Range("I79").Value = 4999
$ws.Range("K79").Value = 4999
$ws.Range("M79").Value = -3907

$ws.Range("H98").Value = 1203.6451
$ws.Range("I98").Value = 1213.7667
$ws.Range("K98").Value = 1213.7667
$ws.Range("M98").Value = 284.2333000000001

$ws.Range("H113").Value = 5652.3687
$ws.Range("I113").Value = 5116.3335
$ws.Range("K113").Value = 5116.3335
$ws.Range("M113").Value = -1862.3335

$ws.Range("H122").Value = 1203.6451
$ws.Range("I122").Value = 1213.7667
$ws.Range("K122").Value = 3641.300099999999
$ws.Range("M122").Value = -1191.300099999999

$ws.Range("H132").Value = 20410160
$ws.Range("I132").Value = 20835340
$ws.Range("K132").Value = 62506020
$ws.Range("M132").Value = -62503490

$ws.Range("H137").Value = 2081.4924
$ws.Range("I137").Value = 1985.4791
$ws.Range("J137").Value = 2324.0527
$ws.Range("K137").Value = 5956.4373
$ws.Range("L137").Value = 6972.158100000001
$ws.Range("M137").Value = -3406.4373
$ws.Range("N137").Value = -12072.1581

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3198
$ws.Range("I2").Value = 3497.75
$ws.Range("K2").Value = 3497.75
$ws.Range("M2").Value = -3384.75

$ws.Range("H32").Value = 6793.7686
$ws.Range("I32").Value = 6268.231
$ws.Range("K32").Value = 6268.231
$ws.Range("M32").Value = -5981.231

$ws.Range("H45").Value = 2594.366
$ws.Range("I45").Value = 2377.48
$ws.Range("K45").Value = 2377.48
$ws.Range("M45").Value = -2000.48

$ws.Range("H74").Value = 3073.9644
$ws.Range("I74").Value = 1429.9412
$ws.Range("K74").Value = 1429.9412
$ws.Range("M74").Value = -555.9412

$ws.Range("H77").Value = 3073.9644
$ws.Range("I77").Value = 1429.9412
$ws.Range("K77").Value = 7149.706
$ws.Range("M77").Value = -2781.706

$ws.Range("H116").Value = 3198
$ws.Range("I116").Value = 3497.75
$ws.Range("K116").Value = 3497.75
$ws.Range("M116").Value = -1203.75

$ws.Range("H132").Value = 1842.6394
$ws.Range("I132").Value = 1785.6666
$ws.Range("K132").Value = 5356.9998
$ws.Range("M132").Value = -2826.9998

$ws.Range("H133").Value = 107170.664
$ws.Range("J133").Value = 107170.664
$ws.Range("L133").Value = 107170.664
$ws.Range("N133").Value = -112230.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3198
$ws.Range("I3").Value = 3497.75
$ws.Range("K3").Value = 3497.75
$ws.Range("M3").Value = -3383.75

$ws.Range("H100").Value = 19250
$ws.Range("J100").Value = 19250
$ws.Range("L100").Value = 19250
$ws.Range("N100").Value = -21414

$ws.Range("H105").Value = 3911.8076
$ws.Range("I105").Value = 3253.0833
$ws.Range("K105").Value = 3253.0833
$ws.Range("M105").Value = -1506.0833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 54813
$ws.Range("J50").Value = 59999.5
$ws.Range("L50").Value = 59999.5
$ws.Range("N50").Value = -61249.5

$ws.Range("H55").Value = 13055.167
$ws.Range("J55").Value = 13055.167
$ws.Range("L55").Value = 13055.167
$ws.Range("N55").Value = -13685.167

$ws.Range("H95").Value = 51541.332
$ws.Range("J95").Value = 51541.332
$ws.Range("L95").Value = 51541.332
$ws.Range("N95").Value = -57033.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 7459.75
$ws.Range("I59").Value = 7419
$ws.Range("K59").Value = 22257
$ws.Range("M59").Value = -21717

$ws.Range("H61").Value = 4200.4165
$ws.Range("I61").Value = 44.5
$ws.Range("J61").Value = 12512.25
$ws.Range("K61").Value = 133.5
$ws.Range("L61").Value = 37536.75
$ws.Range("M61").Value = 81.5
$ws.Range("N61").Value = -37966.75

$ws.Range("H121").Value = 2311.4119
$ws.Range("I121").Value = 1076.909
$ws.Range("J121").Value = 4574.6665
$ws.Range("K121").Value = 3230.727
$ws.Range("L121").Value = 13723.9995
$ws.Range("M121").Value = -1920.727
$ws.Range("N121").Value = -16343.9995

$ws.Range("H124").Value = 5654.6665
$ws.Range("I124").Value = 4779
$ws.Range("K124").Value = 14337
$ws.Range("M124").Value = -9427

$ws.Range("H131").Value = 117882.664
$ws.Range("I131").Value = 715080
$ws.Range("K131").Value = 2145240
$ws.Range("M131").Value = -2140200

$ws.Range("H134").Value = 4142.7144
$ws.Range("I134").Value = 4166.3335
$ws.Range("K134").Value = 12499.0005
$ws.Range("M134").Value = -7429.000499999998

$ws.Range("H136").Value = 8000
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws.Range("H139").Value = 10005303
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

$ws.Range("H140").Value = 4258.857
$ws.Range("I140").Value = 3250.6
$ws.Range("K140").Value = 9751.799999999999
$ws.Range("M140").Value = -4571.799999999999

$ws.Range("H141").Value = 23124.666
$ws.Range("I141").Value = 20949.6
$ws.Range("J141").Value = 34000
$ws.Range("K141").Value = 62848.8
$ws.Range("L141").Value = 102000
$ws.Range("M141").Value = -57668.8
$ws.Range("N141").Value = -112360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 80000
$ws.Range("J62").Value = 80000
$ws.Range("L62").Value = 80000
$ws.Range("N62").Value = -81372

$ws.Range("H65").Value = 80000
$ws.Range("J65").Value = 80000
$ws.Range("L65").Value = 240000
$ws.Range("N65").Value = -246864

$ws.Range("H80").Value = 3709.6956
$ws.Range("I80").Value = 3236.9092
$ws.Range("J80").Value = 4143.0835
$ws.Range("K80").Value = 3236.9092
$ws.Range("L80").Value = 4143.0835
$ws.Range("M80").Value = -2238.9092
$ws.Range("N80").Value = -6139.0835

$ws.Range("H83").Value = 3709.6956
$ws.Range("I83").Value = 3236.9092
$ws.Range("J83").Value = 4143.0835
$ws.Range("K83").Value = 16184.546
$ws.Range("L83").Value = 20715.4175
$ws.Range("M83").Value = -11192.546
$ws.Range("N83").Value = -30699.4175

$ws.Range("H132").Value = 7578640.5
$ws.Range("I132").Value = 9526713
$ws.Range("J132").Value = 2802.4443
$ws.Range("K132").Value = 28580139
$ws.Range("L132").Value = 8407.332900000001
$ws.Range("M132").Value = -28577609
$ws.Range("N132").Value = -13467.3329

$ws.Range("H136").Value = 67457.09
$ws.Range("J136").Value = 67457.09
$ws.Range("L136").Value = 202371.27
$ws.Range("N136").Value = -207471.27

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 42220.5
$ws.Range("I48").Value = 39997
$ws.Range("J48").Value = 44444
$ws.Range("K48").Value = 39997
$ws.Range("L48").Value = 44444
$ws.Range("M48").Value = -39336
$ws.Range("N48").Value = -45766

$ws.Range("H63").Value = 90971.5
$ws.Range("J63").Value = 90971.5
$ws.Range("L63").Value = 90971.5
$ws.Range("N63").Value = -92469.5

$ws.Range("H66").Value = 90971.5
$ws.Range("J66").Value = 90971.5
$ws.Range("L66").Value = 272914.5
$ws.Range("N66").Value = -280402.5

$ws.Range("H92").Value = 99280
$ws.Range("J92").Value = 99280
$ws.Range("L92").Value = 99280
$ws.Range("N92").Value = -104272

$ws.Range("H122").Value = 4222.7607
$ws.Range("I122").Value = 3673.423
$ws.Range("K122").Value = 11020.269
$ws.Range("M122").Value = -8570.269

$ws.Range("H132").Value = 4198.9
$ws.Range("I132").Value = 4197
$ws.Range("J132").Value = 4200.353
$ws.Range("K132").Value = 12591
$ws.Range("L132").Value = 12601.059
$ws.Range("M132").Value = -10061
$ws.Range("N132").Value = -17661.059

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1287.6207
$ws.Range("J100").Value = 1295.3334
$ws.Range("L100").Value = 2590.6668
$ws.Range("N100").Value = -3672.6668

$ws.Range("H122").Value = 2245.7144
$ws.Range("I122").Value = 2245.7144
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6737.1432
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4287.1432
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 1301.2046
$ws.Range("I132").Value = 1199.8235
$ws.Range("K132").Value = 3599.4705
$ws.Range("M132").Value = -1069.4705

$ws.Range("H137").Value = 76618.164
$ws.Range("J137").Value = 76618.164
$ws.Range("L137").Value = 76618.164
$ws.Range("N137").Value = -86818.164
